$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the source feed as the most
# recent record for this market/variety. In the sheet it lands as a new
# row 129, pushing the previously-existing rows 129-156 down to 130-157.
$ws.Rows.Item(129).Insert()

$ws.Range("A129").Value = 7
$ws.Range("B129").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C129").Value = "Ñuble"
$ws.Range("D129").Value = 44543
$ws.Range("E129").Value = 16
$ws.Range("F129").Value = 100112032
$ws.Range("G129").Value = "Zapallo italiano"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 200
$ws.Range("K129").Value = 6000
$ws.Range("L129").Value = 7000
$ws.Range("M129").Value = 6500
$ws.Range("N129").Value = "$/caja 60 unidades"
$ws.Range("O129").Value = "Región del Maule"
$ws.Range("P129").Value = 108
$ws.Range("Q129").Value = 60
$ws.Range("R129").Value = "Hortaliza"
